# Update TPM-derived NATMI edge statistics (ligand/receptor expression, specificity, and edge weights)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.73885133333334
$ws.Range("H2").Value = 107.216554
$ws.Range("I2").Value = 0.01949729408921566
$ws.Range("J2").Value = 0.01949729408921566
$ws.Range("M2").Value = 2.33201
$ws.Range("N2").Value = 6.99603
$ws.Range("O2").Value = 0.3303263034789547
$ws.Range("P2").Value = 0.3303263034789548
$ws.Range("Q2").Value = 83.34335869784667
$ws.Range("R2").Value = 750.0902282806201
$ws.Range("S2").Value = 0.006440469084332682
$ws.Range("T2").Value = 0.006440469084332682

$ws.Range("G3").Value = 35.73885133333334
$ws.Range("H3").Value = 107.216554
$ws.Range("I3").Value = 0.01949729408921566
$ws.Range("J3").Value = 0.01949729408921566
$ws.Range("N3").Value = 5.238131999999999
$ws.Range("O3").Value = 0.2473249515360603
$ws.Range("P3").Value = 0.2473249515360603
$ws.Range("Q3").Value = 62.40160693745866
$ws.Range("R3").Value = 561.6144624371279
$ws.Range("S3").Value = 0.004822167315699578
$ws.Range("T3").Value = 0.004822167315699577

$ws.Range("G4").Value = 35.73885133333334
$ws.Range("H4").Value = 107.216554
$ws.Range("I4").Value = 0.01949729408921566
$ws.Range("J4").Value = 0.01949729408921566
$ws.Range("M4").Value = 1.145780666666667
$ws.Range("N4").Value = 3.437342
$ws.Range("O4").Value = 0.1622984001859565
$ws.Range("P4").Value = 0.1622984001859565
$ws.Range("Q4").Value = 40.94888490660755
$ws.Range("R4").Value = 368.539964159468
$ws.Range("S4").Value = 0.003164379638634807
$ws.Range("T4").Value = 0.003164379638634807

$ws.Range("G5").Value = 35.73885133333334
$ws.Range("H5").Value = 107.216554
$ws.Range("I5").Value = 0.01949729408921566
$ws.Range("J5").Value = 0.01949729408921566
$ws.Range("M5").Value = 1.835881666666667
$ws.Range("N5").Value = 5.507645
$ws.Range("O5").Value = 0.2600503447990285
$ws.Range("P5").Value = 0.2600503447990286
$ws.Range("Q5").Value = 65.61230195059223
$ws.Range("R5").Value = 590.5107175553301
$ws.Range("S5").Value = 0.005070278050548593
$ws.Range("T5").Value = 0.005070278050548593

$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.9215900675332435
$ws.Range("J6").Value = 0.9215900675332435
$ws.Range("M6").Value = 2.33201
$ws.Range("N6").Value = 6.99603
$ws.Range("O6").Value = 0.3303263034789547
$ws.Range("P6").Value = 0.3303263034789548
$ws.Range("Q6").Value = 3939.439556039733
$ws.Range("R6").Value = 35454.9560043576
$ws.Range("S6").Value = 0.3044254403311765
$ws.Range("T6").Value = 0.3044254403311766

$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.9215900675332435
$ws.Range("J7").Value = 0.9215900675332435
$ws.Range("N7").Value = 5.238131999999999
$ws.Range("O7").Value = 0.2473249515360603
$ws.Range("P7").Value = 0.2473249515360603
$ws.Range("Q7").Value = 2949.573458169493
$ws.Range("R7").Value = 26546.16112352543
$ws.Range("S7").Value = 0.227932218788774
$ws.Range("T7").Value = 0.227932218788774

$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.9215900675332435
$ws.Range("J8").Value = 0.9215900675332435
$ws.Range("M8").Value = 1.145780666666667
$ws.Range("N8").Value = 3.437342
$ws.Range("O8").Value = 0.1622984001859565
$ws.Range("P8").Value = 0.1622984001859565
$ws.Range("Q8").Value = 1935.555027985404
$ws.Range("R8").Value = 17419.99525186864
$ws.Range("S8").Value = 0.149572593587913
$ws.Range("T8").Value = 0.149572593587913

$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.9215900675332435
$ws.Range("J9").Value = 0.9215900675332435
$ws.Range("M9").Value = 1.835881666666667
$ws.Range("N9").Value = 5.507645
$ws.Range("O9").Value = 0.2600503447990285
$ws.Range("P9").Value = 0.2600503447990286
$ws.Range("Q9").Value = 3101.335267805377
$ws.Range("R9").Value = 27912.0174102484
$ws.Range("S9").Value = 0.2396598148253799
$ws.Range("T9").Value = 0.23965981482538

$ws.Range("G10").Value = 93.641553
$ws.Range("H10").Value = 280.924659
$ws.Range("I10").Value = 0.05108605424341119
$ws.Range("J10").Value = 0.05108605424341119
$ws.Range("M10").Value = 2.33201
$ws.Range("N10").Value = 6.99603
$ws.Range("O10").Value = 0.3303263034789547
$ws.Range("P10").Value = 0.3303263034789548
$ws.Range("Q10").Value = 218.37303801153
$ws.Range("R10").Value = 1965.35734210377
$ws.Range("S10").Value = 0.01687506745755139
$ws.Range("T10").Value = 0.01687506745755139

$ws.Range("G11").Value = 93.641553
$ws.Range("H11").Value = 280.924659
$ws.Range("I11").Value = 0.05108605424341119
$ws.Range("J11").Value = 0.05108605424341119
$ws.Range("N11").Value = 5.238131999999999
$ws.Range("O11").Value = 0.2473249515360603
$ws.Range("P11").Value = 0.2473249515360603
$ws.Range("Q11").Value = 163.502271766332
$ws.Range("R11").Value = 1471.520445896988
$ws.Range("S11").Value = 0.01263485588992022
$ws.Range("T11").Value = 0.01263485588992022

$ws.Range("G12").Value = 93.641553
$ws.Range("H12").Value = 280.924659
$ws.Range("I12").Value = 0.05108605424341119
$ws.Range("J12").Value = 0.05108605424341119
$ws.Range("M12").Value = 1.145780666666667
$ws.Range("N12").Value = 3.437342
$ws.Range("O12").Value = 0.1622984001859565
$ws.Range("P12").Value = 0.1622984001859565
$ws.Range("Q12").Value = 107.292681024042
$ws.Range("R12").Value = 965.634129216378
$ws.Range("S12").Value = 0.008291184875518628
$ws.Range("T12").Value = 0.008291184875518628

$ws.Range("G13").Value = 93.641553
$ws.Range("H13").Value = 280.924659
$ws.Range("I13").Value = 0.05108605424341119
$ws.Range("J13").Value = 0.05108605424341119
$ws.Range("M13").Value = 1.835881666666667
$ws.Range("N13").Value = 5.507645
$ws.Range("O13").Value = 0.2600503447990285
$ws.Range("P13").Value = 0.2600503447990286
$ws.Range("Q13").Value = 171.914810390895
$ws.Range("R13").Value = 1547.233293518055
$ws.Range("S13").Value = 0.01328494602042095
$ws.Range("T13").Value = 0.01328494602042096

$ws.Range("G14").Value = 14.34625366666667
$ws.Range("H14").Value = 43.038761
$ws.Range("I14").Value = 0.007826584134129748
$ws.Range("J14").Value = 0.007826584134129748
$ws.Range("M14").Value = 2.33201
$ws.Range("N14").Value = 6.99603
$ws.Range("O14").Value = 0.3303263034789547
$ws.Range("P14").Value = 0.3303263034789548
$ws.Range("Q14").Value = 33.45560701320333
$ws.Range("R14").Value = 301.10046311883
$ws.Range("S14").Value = 0.002585326605894115
$ws.Range("T14").Value = 0.002585326605894115

$ws.Range("G15").Value = 14.34625366666667
$ws.Range("H15").Value = 43.038761
$ws.Range("I15").Value = 0.007826584134129748
$ws.Range("J15").Value = 0.007826584134129748
$ws.Range("N15").Value = 5.238131999999999
$ws.Range("O15").Value = 0.2473249515360603
$ws.Range("P15").Value = 0.2473249515360603
$ws.Range("Q15").Value = 25.04919013716133
$ws.Range("R15").Value = 225.442711234452
$ws.Range("S15").Value = 0.001935709541666538
$ws.Range("T15").Value = 0.001935709541666538

$ws.Range("G16").Value = 14.34625366666667
$ws.Range("H16").Value = 43.038761
$ws.Range("I16").Value = 0.007826584134129748
$ws.Range("J16").Value = 0.007826584134129748
$ws.Range("M16").Value = 1.145780666666667
$ws.Range("N16").Value = 3.437342
$ws.Range("O16").Value = 0.1622984001859565
$ws.Range("P16").Value = 0.1622984001859565
$ws.Range("Q16").Value = 16.43766009036245
$ws.Range("R16").Value = 147.938940813262
$ws.Range("S16").Value = 0.001270242083890047
$ws.Range("T16").Value = 0.001270242083890048

$ws.Range("G17").Value = 14.34625366666667
$ws.Range("H17").Value = 43.038761
$ws.Range("I17").Value = 0.007826584134129748
$ws.Range("J17").Value = 0.007826584134129748
$ws.Range("M17").Value = 1.835881666666667
$ws.Range("N17").Value = 5.507645
$ws.Range("O17").Value = 0.2600503447990285
$ws.Range("P17").Value = 0.2600503447990286
$ws.Range("Q17").Value = 26.33802409198278
$ws.Range("R17").Value = 237.042216827845
$ws.Range("S17").Value = 0.002035305902679047
$ws.Range("T17").Value = 0.002035305902679047
